# Add two new German/English vocabulary entries to the "vocab" sheet.
#
# The sheet is sorted by lesson date/number; the two new rows belong to
# lesson 10 (date 2022-04-06 / serial 44657) and are inserted in the
# middle of that block, pushing the existing lesson-10 rows down:
#   - a new row at (current) row 268: "anwesend / abwesend" / "present / absent"
#   - a new row at (current) row 275: "den Zugang verwehren" / "to deny access"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two blank rows first (top one, then the one further down --
# inserting top-down means the second insertion's target index already
# accounts for the shift caused by the first insertion).
$ws.Rows.Item(268).Insert()
$ws.Rows.Item(275).Insert()

# Fill in the lower row.
$ws.Cells.Item(275, 1).Value = "den Zugang verwehren"
$ws.Cells.Item(275, 2).Value = "to deny access"
$ws.Cells.Item(275, 3).Value = 44657
$ws.Cells.Item(275, 4).Value = 10
$ws.Cells.Item(275, 5).Value = "word"

# Fill in the upper row.
$ws.Cells.Item(268, 1).Value = "anwesend / abwesend"
$ws.Cells.Item(268, 2).Value = "present / absent"
$ws.Cells.Item(268, 3).Value = 44657
$ws.Cells.Item(268, 4).Value = 10
$ws.Cells.Item(268, 5).Value = "word"

# Match the saved selection/active cell left behind by the edit.
$ws.Range("A268").Select()
